$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oktober")

# --- Update existing row 14: C14 value changed (101200); D14's running-total formula recalculates ---
$ws.Range("C14").Value = 101200

# --- Row 15 (new transaction: "uang karcis juanda" / Tondo) ---
$ws.Range("A15").Value = 44476
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 13000
$ws.Range("E15").Value = "uang karcis juanda"
$ws.Range("F15").Value = "Tondo"

# --- Row 16 (new transaction: "uang beli bensin pertalite" / yofandi) ---
$ws.Range("A16").Value = 44477
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 10000
$ws.Range("E16").Value = "uang beli bensin pertalite"
$ws.Range("F16").Value = "yofandi"

# --- Extend the running-total formula (D10+B-C pattern) down through the new rows ---
$ws.Range("D11:D16").FormulaR1C1 = "=R[-1]C+RC[-2]-RC[-1]"

# Match the number formatting used by the rest of the running-total column
$ws.Range("D14").Copy()
$ws.Range("D15:D16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Restore active selection to C11 like in the saved file ---
$ws.Range("C11").Select()
